# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/bordered/centered header formatting used by the rest of row 1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill every player's row (2-52) with the team's season record
$lastRow = 52
$ws.Range("AD2:AD" + $lastRow).Value = 55
$ws.Range("AE2:AE" + $lastRow).Value = 107
$ws.Range("AF2:AF" + $lastRow).Value = 0

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
